$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2-3: account holder name / card number
$ws.Range("C2").Value = "Hartmut"
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("C3").Value = "Mohaupt"

# Row 5: opening balance label
$ws.Range("D5").Value = "KONTOSTAND AM 14.06.2025"

# Row 6
$ws.Range("B6").Value = "16.06."
$ws.Range("C6").Value = "17.06."
$ws.Range("D6").Value = "MCDONALDS Dachau"
$ws.Range("E6").Value = "16,49-"

# Row 7
$ws.Range("B7").Value = "17.06."
$ws.Range("C7").Value = "18.06."
$ws.Range("D7").Value = "KARTENZAHLUNG ARAL TANKSTELLE"
$ws.Range("E7").Value = "71,58-"

# Row 8
$ws.Range("B8").Value = "20.06."
$ws.Range("C8").Value = "21.06."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 11881583"
$ws.Range("E8").Value = "87,77-"

# Row 9
$ws.Range("B9").Value = "23.06."
$ws.Range("C9").Value = "24.06."
$ws.Range("D9").Value = "MITGLIEDSBEITRAG ZEUS BODYPOWER"
$ws.Range("E9").Value = "24,52-"

# Row 10
$ws.Range("B10").Value = "26.06."
$ws.Range("C10").Value = "27.06."
$ws.Range("D10").Value = "KARTENZAHLUNG JET TANKSTELLE"
$ws.Range("E10").Value = "64,83-"

# Row 11
$ws.Range("B11").Value = "29.06."
$ws.Range("C11").Value = "30.06."
$ws.Range("D11").Value = "PAYPAL RFDWBF"
$ws.Range("E11").Value = "64,79-"

# Row 12: closing balance
$ws.Range("D12").Value = "KONTOSTAND AM 02.07.2025"
$ws.Range("E12").Value = "329,98-"

# Row 13: next statement date
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 10.07.2025"
